$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("SWR")
$ws2 = $wb.Worksheets.Item("settings")

# --- Rename the recipe/designator labels to uppercase (core change) ---
$ws1.Range("A2").Clear()
$ws1.Range("A2").Value = "123456-EXTRAPLACE"

$ws1.Range("A3").Clear()
$ws1.Range("A3").Value = "123457-NOPLACE-ALL"

$ws1.Range("A4").Clear()
$ws1.Range("A4").Value = "123458-NOPLACE-PARTIAL"

$ws1.Range("A5").Clear()
$ws1.Range("A5").Value = "123459-PARTSUB-ALL"

$ws1.Range("A6").Clear()
$ws1.Range("A6").Value = "123460-PARTSUB-PARTIAL"

# --- Leftover UPPER() helper formulas (rows 18-21, col C) ---
$ws1.Range("C18").Formula = "=UPPER(A8)"
$ws1.Range("C19").Formula = "=UPPER(A9)"
$ws1.Range("C20").Formula = "=UPPER(A10)"
$ws1.Range("C21").Formula = "=UPPER(A11)"

# --- Drop the now-unused outline levels on SWR sheet ---
$ws1.Rows.Ungroup()

# --- Unfreeze panes / reset selection on SWR sheet ---
$ws1.Activate()
$excel.ActiveWindow.FreezePanes = $false
$ws1.Range("F13").Select()

# --- Narrow duplicate-values conditional formatting away from the renamed rows ---
$cf = $ws1.Range("A1").FormatConditions.Item(1)
$cf.ModifyAppliesToRange($ws1.Range("A7:A1048576"))

# --- Unfreeze panes / reset selection on settings sheet ---
$ws2.Activate()
$excel.ActiveWindow.FreezePanes = $false
$ws2.Range("B6").Select()

# Leave SWR as the active/selected sheet, matching the original workbook state
$ws1.Activate()
